$wb = $excel.ActiveWorkbook

# Suppress the "this sheet contains data" confirmation Excel normally raises
# when deleting a non-empty worksheet.
$excel.DisplayAlerts = $false

# Remove the "Correct order for AI" sheet and keep only the main data sheet,
# renamed to "Sheet1" (matches the workbook's new single-sheet layout).
$wsRemove = $wb.Worksheets.Item("Correct order for AI")
$wsRemove.Delete()

$ws = $wb.Worksheets.Item("Across Survyes")
$ws.Name = "Sheet1"
$ws.Select()

$excel.DisplayAlerts = $true
